$wb = $excel.ActiveWorkbook

# Update the "lighting" row's share value on the SoCEUtiNTY sheet: replace the
# formula-derived value in B6 with a hard-coded literal. Dependent formulas in
# C6 (=B6) and D6 (shared formula) recalc automatically.
$soc = $wb.Worksheets("SoCEUtiNTY")
$soc.Range("B6").Value = 0.08

# Move the active tab / selection from "About" to "SoCEUtiNTY", and update the
# remembered selection on SoCEUtiNTY from B8 to B7.
$soc.Activate()
[void]$soc.Range("B7").Select()

# Turn off Excel's multithreaded/concurrent calculation setting.
$excel.Application.MultiThreadedCalculation.Enabled = $false
